$wb = $excel.ActiveWorkbook

$wsMoCo = $wb.Worksheets.Item("MoCo")
$wsClf  = $wb.Worksheets.Item("MoCoClf")

# ---------------------------------------------------------------------------
# 1) MoCo sheet (sheet2): append new experiment row 37
# ---------------------------------------------------------------------------
$wsMoCo.Cells.Item(37, 1).Value = 8852525
$wsMoCo.Cells.Item(37, 2).Value = "OrdLabelMoCo (alpha=2, ratio, start from 8782858_480)"
$wsMoCo.Cells.Item(37, 3).Value = 3
$wsMoCo.Cells.Item(37, 4).Value = "SGD"
$wsMoCo.Cells.Item(37, 5).Value = 128
$wsMoCo.Cells.Item(37, 6).Value = 1
$wsMoCo.Cells.Item(37, 7).Value = 500

# ---------------------------------------------------------------------------
# 2) MoCoClf sheet (sheet3): fill in the validation-accuracy column (K) for
#    the existing rows 35-37, then append four new experiment rows (38-41).
# ---------------------------------------------------------------------------
$wsClf.Cells.Item(35, 11).Value = 0.6967
$wsClf.Cells.Item(35, 11).NumberFormat = "0.00%"

$wsClf.Cells.Item(36, 11).Value = 0.8006
$wsClf.Cells.Item(36, 11).NumberFormat = "0.00%"
$wsClf.Cells.Item(36, 11).Interior.Color = 49407

$wsClf.Cells.Item(37, 11).Value = 0.6303
$wsClf.Cells.Item(37, 11).NumberFormat = "0.00%"
$wsClf.Cells.Item(37, 11).Interior.Color = 49407

# copy the formatting of row 37 onto the four freshly-appended rows first so
# the new cells pick up the same style (fill / borders) as their neighbours
$wsClf.Range("A37:I37").Copy()
$wsClf.Range("A38:I41").PasteSpecial(-4122)  # xlPasteFormats

$wsClf.Cells.Item(38, 1).Value = 8855877
$wsClf.Cells.Item(38, 2).Value = "MoCoCLfV2"
$wsClf.Cells.Item(38, 3).Value = 0.01
$wsClf.Cells.Item(38, 4).Value = "SGD"
$wsClf.Cells.Item(38, 5).Value = 128
$wsClf.Cells.Item(38, 6).Value = 1
$wsClf.Cells.Item(38, 7).Value = 200
$wsClf.Cells.Item(38, 8).Value = "8852525_580"
$wsClf.Cells.Item(38, 9).Value = ""

$wsClf.Cells.Item(39, 1).Value = 8855878
$wsClf.Cells.Item(39, 2).Value = "MoCoCLfV2"
$wsClf.Cells.Item(39, 3).Value = 0.01
$wsClf.Cells.Item(39, 4).Value = "SGD"
$wsClf.Cells.Item(39, 5).Value = 128
$wsClf.Cells.Item(39, 6).Value = 1
$wsClf.Cells.Item(39, 7).Value = 200
$wsClf.Cells.Item(39, 8).Value = "8852525_640"
$wsClf.Cells.Item(39, 9).Value = ""

$wsClf.Cells.Item(40, 1).Value = 8855876
$wsClf.Cells.Item(40, 2).Value = "MoCoClfV2Fea + n lanes + speed"
$wsClf.Cells.Item(40, 3).Value = 0.01
$wsClf.Cells.Item(40, 4).Value = "SGD"
$wsClf.Cells.Item(40, 5).Value = 128
$wsClf.Cells.Item(40, 6).Value = 1
$wsClf.Cells.Item(40, 7).Value = 200
$wsClf.Cells.Item(40, 8).Value = "8852525_580"
$wsClf.Cells.Item(40, 9).Value = ""

$wsClf.Cells.Item(41, 1).Value = 8855875
$wsClf.Cells.Item(41, 2).Value = "MoCoClfV2Fea + n lanes + speed"
$wsClf.Cells.Item(41, 3).Value = 0.01
$wsClf.Cells.Item(41, 4).Value = "SGD"
$wsClf.Cells.Item(41, 5).Value = 128
$wsClf.Cells.Item(41, 6).Value = 1
$wsClf.Cells.Item(41, 7).Value = 200
$wsClf.Cells.Item(41, 8).Value = "8852525_640"
$wsClf.Cells.Item(41, 9).Value = ""

# ---------------------------------------------------------------------------
# 3) View-state: MoCoClf becomes the active/selected sheet, each sheet keeps
#    its own scroll position / active-cell selection.
# ---------------------------------------------------------------------------
$wsMoCo.Select()
$wsMoCo.Range("B40").Select()

$wsClf.Select()
$wsClf.Range("A40").Select()
